$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: change status from "Pendente" to "Concluído" and restyle the row
# to match the "completed" look used elsewhere (fill highlight style, same as
# row 4/9/12 etc.) while keeping text/values otherwise intact.
$ws.Range("A9:F9").Copy() | Out-Null
$ws.Range("A11:F11").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").Value = "Concluído"

# --- Row 14: restyle to the same highlighted look (values/text unchanged).
$ws.Range("A9:F9").Copy() | Out-Null
$ws.Range("A14:F14").PasteSpecial(-4122) | Out-Null

# --- Row 15: fill in the new backlog item (previously blank placeholder row).
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Sistema Completo"
$ws.Range("C15").Value = "Criar Gameficação"
$ws.Range("D15").Value = "Pendente"
$ws.Range("E15").Value = "Normal"

# --- Update the active selection to match where the user ended up (F15).
$ws.Range("F15").Select() | Out-Null
